# Anapa, GK (Goryachiy Klyuch), Gelendzhik. 2018 (added)
#
# This script reproduces, against the ORIGINAL (pre-edit) row numbering of
# sheet1, the following logical changes:
#   1. "Goryachiy Klyuch" 2018 row (row 10) gains previously-missing values.
#   2. A new "Gelendzhik" 2017 stub row is inserted right before the
#      existing (until now empty) "Gelendzhik" 2018 row (row 22), which
#      itself gets filled in with the full 2018 data.
#   3. A new "Anapa" 2018 data row is inserted right after the existing
#      "Anapa" row (row 16), whose own year is corrected 2018 -> 2017 and
#      which gains a saldo (U) value. Doing this insert last keeps the
#      Gelendzhik row numbers above stable while writing them.
#
# xlCenter = -4108 (matches the workbook's existing direct cell style that
# centers numbers/text -- reusing it keeps the same style index instead of
# Excel minting a new one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

function Set-Cell($addr, $value) {
    $ws.Range($addr).HorizontalAlignment = $xlCenter
    $ws.Range($addr).Value = $value
}

# ---------------------------------------------------------------------
# 1. Goryachiy Klyuch, 2018 (row 10) - fill previously blank columns
# ---------------------------------------------------------------------
Set-Cell "D10" 8.8
Set-Cell "E10" 216
Set-Cell "N10" 996
Set-Cell "Q10" 72.433
Set-Cell "R10" 982

# ---------------------------------------------------------------------
# 2. Gelendzhik: insert the new 2017 row ahead of the 2018 row, then
#    populate both (do this BEFORE the Anapa insert below so the row
#    numbers used here still refer to the original layout).
# ---------------------------------------------------------------------
$ws.Rows.Item(23).Insert()

# Row 22 was the "Gelendzhik 2018" stub (name + year only); it becomes
# the new "Gelendzhik 2017" row.
Set-Cell "B22" 2017
Set-Cell "U22" 2138

# Row 23 is the freshly inserted blank row; it becomes the full
# "Gelendzhik 2018" data row.
Set-Cell "A23" "Геленджик"
Set-Cell "B23" 2018
Set-Cell "D23" 17
Set-Cell "E23" 196
Set-Cell "F23" 33853
Set-Cell "N23" 3047
Set-Cell "O23" 1997.1
Set-Cell "P23" 174.9
Set-Cell "Q23" 111.617
Set-Cell "R23" 1622
Set-Cell "S23" 13072.3
Set-Cell "T23" 753.8
Set-Cell "U23" 160

# ---------------------------------------------------------------------
# 3. Anapa: correct row 16 to 2017 and insert the new 2018 row after it.
# ---------------------------------------------------------------------
Set-Cell "B16" 2017
Set-Cell "U16" -677

$ws.Rows.Item(17).Insert()

Set-Cell "A17" "Анапа"
Set-Cell "B17" 2018
Set-Cell "D17" 24.9
Set-Cell "E17" 507
Set-Cell "F17" 31629
Set-Cell "N17" 4471
Set-Cell "O17" 1435.5
Set-Cell "P17" 4022.4
Set-Cell "Q17" 422.519
Set-Cell "R17" 7064
Set-Cell "S17" 20304.5
Set-Cell "T17" 389.4
Set-Cell "U17" 7470

# ---------------------------------------------------------------------
# Restore the selection to roughly where the author last left off.
# ---------------------------------------------------------------------
$ws.Range("J34").Select()
